# Fix from Yizhen for \omic\ in cds corrected test cases in library selection cds
# Replace the RANDOM "Files" stat query in cell C2 of the "startup" sheet with
# the new RANDOM "Participants" query (using apoc.coll.sort).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newQuery = "Match (f)<--(g:genomic_info)`nWHERE g.library_selection in ['RANDOM']`nMATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)`nWITH p, s, apoc.coll.sort(collect(distinct(samp.sample_id))) as samples`nRETURN `ncoalesce(p.participant_id,'') as ``Participant ID``,`ncoalesce(s.study_name, '') as ``Study Name``,`ncoalesce(s.phs_accession,'') as ``Accession``,`ncoalesce(p.gender,'') as ``Gender``,`ncoalesce(apoc.text.join(samples, ','), '') as ``Samples```nORDER BY ``Participant ID``LIMIT 100"

$ws.Range("C2").Value = $newQuery

# Update the view to match: top-left cell C1 and selection at C2
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("C2").Select()

$wb.Save()
